$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-82 currently hold the date serial 45178
# (2023-09-09) and need to be bumped to 45179 (2023-09-10).
$ws.Range("C2:C82").Value = 45179
